$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) column cells that look like numbers stay as text, matching source formatting
$textCells = @(
    "D2", "D3", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16",
    "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D30",
    "D31", "D32", "D33", "D34", "D35", "D39", "D40", "D41", "D42", "D43", "D44", "D45",
    "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.831.36"
$ws.Range("E2").Value = "  -0.71%  "

$ws.Range("D3").Value = "3.949.80"
$ws.Range("E3").Value = "  +3.93%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "601.76"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").Value = "166.51"
$ws.Range("E6").Value = "  +0.68%  "

$ws.Range("D7").Value = "3.945.32"
$ws.Range("E7").Value = "  +3.89%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  -1.60%  "

$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  -2.96%  "

$ws.Range("D11").Value = "6.37"
$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  +0.62%  "

$ws.Range("D14").Value = "37.24"
$ws.Range("E14").Value = "  -0.71%  "

$ws.Range("D15").Value = "4.593.87"
$ws.Range("E15").Value = "  +3.59%  "

$ws.Range("D16").Value = "3.936.70"
$ws.Range("E16").Value = "  +3.55%  "

$ws.Range("D17").Value = "68.955.28"
$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").Value = "7.44"
$ws.Range("E18").Value = "  -0.58%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "17.13"
$ws.Range("E19").Value = "  -2.38%  "

$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "0.112"
$ws.Range("E20").Value = "  -1.39%  "

$ws.Range("D21").Value = "11.21"
$ws.Range("E21").Value = "  -0.87%  "

$ws.Range("D22").Value = "487.23"
$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("D23").Value = "0.0000172"
$ws.Range("E23").Value = "  +13.52%  "

$ws.Range("D24").Value = "0.724"
$ws.Range("E24").Value = "  -0.27%  "

$ws.Range("D25").Value = "84.78"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("E26").Value = "  -1.13%  "

$ws.Range("D27").Value = "12.04"
$ws.Range("E27").Value = "  -2.23%  "

$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  -0.28%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").Value = "2.94"
$ws.Range("E30").Value = "  -1.78%  "

$ws.Range("D31").Value = "4.087.23"
$ws.Range("E31").Value = "  +3.58%  "

$ws.Range("D32").Value = "2.40"
$ws.Range("E32").Value = "  -0.67%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "32.36"
$ws.Range("E33").Value = "  +0.60%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "7.81"
$ws.Range("E34").Value = "  -4.00%  "

$ws.Range("D35").Value = "3.886.15"
$ws.Range("E35").Value = "  +3.74%  "

$ws.Range("E36").Value = "  -0.67%  "

$ws.Range("E37").Value = "  +2.17%  "

$ws.Range("E38").Value = "  -0.82%  "

$ws.Range("D39").Value = "5.92"
$ws.Range("E39").Value = "  -0.81%  "

$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "3.16"
$ws.Range("E40").Value = "  +3.34%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").Value = "0.318"
$ws.Range("E42").Value = "  -2.15%  "

$ws.Range("D43").Value = "434.89"
$ws.Range("E43").Value = "  +2.23%  "

$ws.Range("D44").Value = "48.48"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").Value = "1.99"
$ws.Range("E45").Value = "  -0.21%  "

$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").Value = "8.51"
$ws.Range("E46").Value = "  +0.77%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").Value = "26.77"
$ws.Range("E48").Value = "  +9.64%  "

$ws.Range("D49").Value = "2.841.17"
$ws.Range("E49").Value = "  +0.72%  "

$ws.Range("D50").Value = "141.87"
$ws.Range("E50").Value = "  -0.37%  "

$ws.Range("D51").Value = "0.000267"
$ws.Range("E51").Value = "  +18.44%  "
